$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.773.31"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").Value = "2.045.32"
$ws.Range("E3").Value = "  +0.90%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.611"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.11"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.55%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -1.67%  "

$ws.Range("E10").Value = "  +2.47%  "

$ws.Range("E11").Value = "  +0.14%  "

$ws.Range("D12").Value = "2.348.09"
$ws.Range("E12").Value = "  +0.84%  "

$ws.Range("E13").Value = "  -0.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.764"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.04%  "

$ws.Range("D17").Value = "2.047.60"
$ws.Range("E17").Value = "  +1.11%  "

$ws.Range("D18").Value = "37.747.41"
$ws.Range("E18").Value = "  +0.21%  "

$ws.Range("E19").Value = "  -0.61%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.83%  "

$ws.Range("E21").Value = "  +0.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "222.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.04%  "

$ws.Range("E23").Value = "  +0.51%  "

$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("E25").Value = "  +3.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.46%  "

$ws.Range("E27").Value = "  +0.99%  "

$ws.Range("E28").Value = "  +0.18%  "

$ws.Range("E29").Value = "  -0.76%  "

$ws.Range("E30").Value = "  +0.22%  "

$ws.Range("E31").Value = "  -0.49%  "

$ws.Range("E32").Value = "  +8.06%  "

$ws.Range("E33").Value = "  -0.87%  "

$ws.Range("E34").Value = "  +0.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0602"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.65%  "

$ws.Range("E37").Value = "  +4.56%  "

$ws.Range("E38").Value = "  +7.83%  "

$ws.Range("E39").Value = "  +0.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.12%  "

$ws.Range("D41").Value = "1.531.11"
$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.04%  "

$ws.Range("E43").Value = "  -0.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0887"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.08%  "

$ws.Range("E47").Value = "  +0.35%  "

$ws.Range("E48").Value = "  +0.28%  "

$ws.Range("E49").Value = "  -0.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.59%  "

$ws.Range("D51").Value = "2.236.87"
$ws.Range("E51").Value = "  +0.89%  "
